# Commit 2019-08-01 kl. 09:24
# Marks more days as "done" on the fix-zone tracker (week of rows 19-20),
# updates the selection, and lets the dependent totals recalc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: Fri (H19) gets marked done (apply the same "done" look used
#     on the rest of that row / the row above), and the done-count (I19)
#     goes from 4 to 5.
$ws.Range("D19").Copy()
$ws.Range("H19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I19").Value = 5

# --- Row 20: Mon/Tue/Wed (D20:F20) get marked done too, copying the
#     "done" formatting already used on row 19, and the done-count (I20)
#     is filled in with 3.
$ws.Range("D19").Copy()
$ws.Range("D20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E19").Copy()
$ws.Range("E20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F19").Copy()
$ws.Range("F20").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I20").Value = 3

# --- Move the selection (no other content change intended here).
[void]$ws.Range("A2").Select()

$excel.CutCopyMode = 0
